# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (G) previously held a "Strike#" style value; this re-derives
# and writes the new strikeout (K) values for each game row on the active
# sheet, replacing the stale numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row (A-index+2) -> new K value, keyed by the worksheet row number
$newK = @{
    2  = 2
    3  = 1
    4  = 1
    5  = 1
    6  = 0
    7  = 1
    8  = 1
    9  = 1
    10 = 2
    11 = 1
    12 = 1
    13 = 1
    14 = 1
    15 = 0
    16 = 2
    17 = 1
    18 = 0
    21 = 1
    24 = 1
    25 = 2
    26 = 2
}

foreach ($r in $newK.Keys) {
    $ws.Cells.Item($r, 7).Value = $newK[$r]
}
